# Nodal-Acvr2a.xlsx update: refresh NATMI TPM numbers and add the
# "Inflammatory-Mac" sending-cluster block (now 4 sending clusters x
# 5 target clusters = 20 data rows instead of 15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20

# row 2: ECs -> ECs
$data[0,0] = 'ECs'
$data[0,1] = 'Nodal'
$data[0,2] = 'Acvr2a'
$data[0,3] = 'ECs'
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 0.7356009999999999
$data[0,7] = 2.206803
$data[0,8] = 0.5568025227146887
$data[0,9] = 0.5688250870390384
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 9.24193
$data[0,13] = 27.72579
$data[0,14] = 0.1468938537243544
$data[0,15] = 0.1569651396557324
$data[0,16] = 6.798372949929999
$data[0,17] = 61.18535654937
$data[0,18] = 0.08179086832500301
$data[0,19] = 0.0892857092267668

# row 3: ECs -> FAPs
$data[1,0] = 'ECs'
$data[1,1] = 'Nodal'
$data[1,2] = 'Acvr2a'
$data[1,3] = 'FAPs'
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.7356009999999999
$data[1,7] = 2.206803
$data[1,8] = 0.5568025227146887
$data[1,9] = 0.5688250870390384
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 29.54200233333333
$data[1,13] = 88.626007
$data[1,14] = 0.469548954544906
$data[1,15] = 0.5017420086455576
$data[1,16] = 21.73112645840233
$data[1,17] = 195.580138125621
$data[1,18] = 0.2614460424286483
$data[1,19] = 0.2854034417389513

# row 4: ECs -> Inflammatory-Mac
$data[2,0] = 'ECs'
$data[2,1] = 'Nodal'
$data[2,2] = 'Acvr2a'
$data[2,3] = 'Inflammatory-Mac'
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.7356009999999999
$data[2,7] = 2.206803
$data[2,8] = 0.5568025227146887
$data[2,9] = 0.5688250870390384
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 7.349831333333333
$data[2,13] = 22.049494
$data[2,14] = 0.1168203014713749
$data[2,15] = 0.1248296948454213
$data[2,16] = 5.406543278631333
$data[2,17] = 48.658889507682
$data[2,18] = 0.06504583856355198
$data[2,19] = 0.07100626203550338

# row 5: ECs -> MuSCs
$data[3,0] = 'ECs'
$data[3,1] = 'Nodal'
$data[3,2] = 'Acvr2a'
$data[3,3] = 'MuSCs'
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.7356009999999999
$data[3,7] = 2.206803
$data[3,8] = 0.5568025227146887
$data[3,9] = 0.5688250870390384
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 12.1104985
$data[3,13] = 24.220997
$data[3,14] = 0.1924876941491673
$data[3,15] = 0.1371233128688515
$data[3,16] = 8.908494807098498
$data[3,17] = 53.45096884259099
$data[3,18] = 0.1071776336937898
$data[3,19] = 0.07799918037770577

# row 6: ECs -> Resolving-Mac
$data[4,0] = 'ECs'
$data[4,1] = 'Nodal'
$data[4,2] = 'Acvr2a'
$data[4,3] = 'Resolving-Mac'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.7356009999999999
$data[4,7] = 2.206803
$data[4,8] = 0.5568025227146887
$data[4,9] = 0.5688250870390384
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 4.671440333333334
$data[4,13] = 14.014321
$data[4,14] = 0.07424919611019735
$data[4,15] = 0.079339843984437
$data[4,16] = 3.436316180640334
$data[4,17] = 30.926845625763
$data[4,18] = 0.04134213970369553
$data[4,19] = 0.04513049366011111

# row 7: FAPs -> ECs
$data[5,0] = 'FAPs'
$data[5,1] = 'Nodal'
$data[5,2] = 'Acvr2a'
$data[5,3] = 'ECs'
$data[5,4] = 2
$data[5,5] = 0.6666666666666666
$data[5,6] = 0.3152516666666667
$data[5,7] = 0.9457549999999999
$data[5,8] = 0.2386251830680085
$data[5,9] = 0.2437776141289484
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 9.24193
$data[5,13] = 27.72579
$data[5,14] = 0.1468938537243544
$data[5,15] = 0.1569651396557324
$data[5,16] = 2.913533835716667
$data[5,17] = 26.22180452145
$data[5,18] = 0.03505257273653934
$data[5,19] = 0.03826458724669164

# row 8: FAPs -> FAPs
$data[6,0] = 'FAPs'
$data[6,1] = 'Nodal'
$data[6,2] = 'Acvr2a'
$data[6,3] = 'FAPs'
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 0.3152516666666667
$data[6,7] = 0.9457549999999999
$data[6,8] = 0.2386251830680085
$data[6,9] = 0.2437776141289484
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 29.54200233333333
$data[6,13] = 88.626007
$data[6,14] = 0.469548954544906
$data[6,15] = 0.5017420086455576
$data[6,16] = 9.313165472253889
$data[6,17] = 83.818489250285
$data[6,18] = 0.1120462052376702
$data[6,19] = 0.1223134697758803

# row 9: FAPs -> Inflammatory-Mac
$data[7,0] = 'FAPs'
$data[7,1] = 'Nodal'
$data[7,2] = 'Acvr2a'
$data[7,3] = 'Inflammatory-Mac'
$data[7,4] = 2
$data[7,5] = 0.6666666666666666
$data[7,6] = 0.3152516666666667
$data[7,7] = 0.9457549999999999
$data[7,8] = 0.2386251830680085
$data[7,9] = 0.2437776141289484
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 7.349831333333333
$data[7,13] = 22.049494
$data[7,14] = 0.1168203014713749
$data[7,15] = 0.1248296948454213
$data[7,16] = 2.317046577552222
$data[7,17] = 20.85341919797
$data[7,18] = 0.02787626582466677
$data[7,19] = 0.0304306851818615

# row 10: FAPs -> MuSCs
$data[8,0] = 'FAPs'
$data[8,1] = 'Nodal'
$data[8,2] = 'Acvr2a'
$data[8,3] = 'MuSCs'
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.3152516666666667
$data[8,7] = 0.9457549999999999
$data[8,8] = 0.2386251830680085
$data[8,9] = 0.2437776141289484
$data[8,10] = 2
$data[8,11] = 1
$data[8,12] = 12.1104985
$data[8,13] = 24.220997
$data[8,14] = 0.1924876941491673
$data[8,15] = 0.1371233128688515
$data[8,16] = 3.817854836289166
$data[8,17] = 22.907129017735
$data[8,18] = 0.04593241125468388
$data[8,19] = 0.03342759405262596

# row 11: FAPs -> Resolving-Mac
$data[9,0] = 'FAPs'
$data[9,1] = 'Nodal'
$data[9,2] = 'Acvr2a'
$data[9,3] = 'Resolving-Mac'
$data[9,4] = 2
$data[9,5] = 0.6666666666666666
$data[9,6] = 0.3152516666666667
$data[9,7] = 0.9457549999999999
$data[9,8] = 0.2386251830680085
$data[9,9] = 0.2437776141289484
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 4.671440333333334
$data[9,13] = 14.014321
$data[9,14] = 0.07424919611019735
$data[9,15] = 0.079339843984437
$data[9,16] = 1.472679350817222
$data[9,17] = 13.254114157355
$data[9,18] = 0.01771772801444831
$data[9,19] = 0.01934127787188905

# row 12: Inflammatory-Mac -> ECs
$data[10,0] = 'Inflammatory-Mac'
$data[10,1] = 'Nodal'
$data[10,2] = 'Acvr2a'
$data[10,3] = 'ECs'
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.1864953333333333
$data[10,7] = 0.559486
$data[10,8] = 0.1411649414213912
$data[10,9] = 0.1442129961972698
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 9.24193
$data[10,13] = 27.72579
$data[10,14] = 0.1468938537243544
$data[10,15] = 0.1569651396557324
$data[10,16] = 1.723576815993333
$data[10,17] = 15.51219134394
$data[10,18] = 0.0207362622561609
$data[10,19] = 0.02263641308827605

# row 13: Inflammatory-Mac -> FAPs
$data[11,0] = 'Inflammatory-Mac'
$data[11,1] = 'Nodal'
$data[11,2] = 'Acvr2a'
$data[11,3] = 'FAPs'
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.1864953333333333
$data[11,7] = 0.559486
$data[11,8] = 0.1411649414213912
$data[11,9] = 0.1442129961972698
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 29.54200233333333
$data[11,13] = 88.626007
$data[11,14] = 0.469548954544906
$data[11,15] = 0.5017420086455576
$data[11,16] = 5.509445572489112
$data[11,17] = 49.58501015240201
$data[11,18] = 0.06628385066280715
$data[11,19] = 0.07235771838481229

# row 14: Inflammatory-Mac -> Inflammatory-Mac
$data[12,0] = 'Inflammatory-Mac'
$data[12,1] = 'Nodal'
$data[12,2] = 'Acvr2a'
$data[12,3] = 'Inflammatory-Mac'
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.1864953333333333
$data[12,7] = 0.559486
$data[12,8] = 0.1411649414213912
$data[12,9] = 0.1442129961972698
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 7.349831333333333
$data[12,13] = 22.049494
$data[12,14] = 0.1168203014713749
$data[12,15] = 0.1248296948454213
$data[12,16] = 1.370709244453778
$data[12,17] = 12.336383200084
$data[12,18] = 0.01649093101403589
$data[12,19] = 0.01800206430804909

# row 15: Inflammatory-Mac -> MuSCs
$data[13,0] = 'Inflammatory-Mac'
$data[13,1] = 'Nodal'
$data[13,2] = 'Acvr2a'
$data[13,3] = 'MuSCs'
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.1864953333333333
$data[13,7] = 0.559486
$data[13,8] = 0.1411649414213912
$data[13,9] = 0.1442129961972698
$data[13,10] = 2
$data[13,11] = 1
$data[13,12] = 12.1104985
$data[13,13] = 24.220997
$data[13,14] = 0.1924876941491673
$data[13,15] = 0.1371233128688515
$data[13,16] = 2.258551454590333
$data[13,17] = 13.551308727542
$data[13,18] = 0.02717251406890587
$data[13,19] = 0.01977496379731272

# row 16: Inflammatory-Mac -> Resolving-Mac
$data[14,0] = 'Inflammatory-Mac'
$data[14,1] = 'Nodal'
$data[14,2] = 'Acvr2a'
$data[14,3] = 'Resolving-Mac'
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.1864953333333333
$data[14,7] = 0.559486
$data[14,8] = 0.1411649414213912
$data[14,9] = 0.1442129961972698
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 4.671440333333334
$data[14,13] = 14.014321
$data[14,14] = 0.07424919611019735
$data[14,15] = 0.079339843984437
$data[14,16] = 0.871201822111778
$data[14,17] = 7.840816399006001
$data[14,18] = 0.0104813834194814
$data[14,19] = 0.01144183661881959

# row 17: MuSCs -> ECs
$data[15,0] = 'MuSCs'
$data[15,1] = 'Nodal'
$data[15,2] = 'Acvr2a'
$data[15,3] = 'ECs'
$data[15,4] = 1
$data[15,5] = 0.5
$data[15,6] = 0.0837685
$data[15,7] = 0.167537
$data[15,8] = 0.06340735279591164
$data[15,9] = 0.04318430263474329
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 9.24193
$data[15,13] = 27.72579
$data[15,14] = 0.1468938537243544
$data[15,15] = 0.1569651396557324
$data[15,16] = 0.7741826132049999
$data[15,17] = 4.64509567923
$data[15,18] = 0.00931415040665118
$data[15,19] = 0.006778430093997892

# row 18: MuSCs -> FAPs
$data[16,0] = 'MuSCs'
$data[16,1] = 'Nodal'
$data[16,2] = 'Acvr2a'
$data[16,3] = 'FAPs'
$data[16,4] = 1
$data[16,5] = 0.5
$data[16,6] = 0.0837685
$data[16,7] = 0.167537
$data[16,8] = 0.06340735279591164
$data[16,9] = 0.04318430263474329
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 29.54200233333333
$data[16,13] = 88.626007
$data[16,14] = 0.469548954544906
$data[16,15] = 0.5017420086455576
$data[16,16] = 2.474689222459833
$data[16,17] = 14.848135334759
$data[16,18] = 0.02977285621578033
$data[16,19] = 0.02166737874591374

# row 19: MuSCs -> Inflammatory-Mac
$data[17,0] = 'MuSCs'
$data[17,1] = 'Nodal'
$data[17,2] = 'Acvr2a'
$data[17,3] = 'Inflammatory-Mac'
$data[17,4] = 1
$data[17,5] = 0.5
$data[17,6] = 0.0837685
$data[17,7] = 0.167537
$data[17,8] = 0.06340735279591164
$data[17,9] = 0.04318430263474329
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 7.349831333333333
$data[17,13] = 22.049494
$data[17,14] = 0.1168203014713749
$data[17,15] = 0.1248296948454213
$data[17,16] = 0.6156843460463333
$data[17,17] = 3.694106076278
$data[17,18] = 0.007407266069120222
$data[17,19] = 0.005390683320007327

# row 20: MuSCs -> MuSCs
$data[18,0] = 'MuSCs'
$data[18,1] = 'Nodal'
$data[18,2] = 'Acvr2a'
$data[18,3] = 'MuSCs'
$data[18,4] = 1
$data[18,5] = 0.5
$data[18,6] = 0.0837685
$data[18,7] = 0.167537
$data[18,8] = 0.06340735279591164
$data[18,9] = 0.04318430263474329
$data[18,10] = 2
$data[18,11] = 1
$data[18,12] = 12.1104985
$data[18,13] = 24.220997
$data[18,14] = 0.1924876941491673
$data[18,15] = 0.1371233128688515
$data[18,16] = 1.01447829359725
$data[18,17] = 4.057913174388999
$data[18,18] = 0.01220513513178779
$data[18,19] = 0.005921574641207072

# row 21: MuSCs -> Resolving-Mac
$data[19,0] = 'MuSCs'
$data[19,1] = 'Nodal'
$data[19,2] = 'Acvr2a'
$data[19,3] = 'Resolving-Mac'
$data[19,4] = 1
$data[19,5] = 0.5
$data[19,6] = 0.0837685
$data[19,7] = 0.167537
$data[19,8] = 0.06340735279591164
$data[19,9] = 0.04318430263474329
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 4.671440333333334
$data[19,13] = 14.014321
$data[19,14] = 0.07424919611019735
$data[19,15] = 0.079339843984437
$data[19,16] = 0.3913195495628334
$data[19,17] = 2.347917297377
$data[19,18] = 0.004707944972572114
$data[19,19] = 0.003426235833617244

# Write the full A2:T21 block in one shot and refresh the used range.
$ws.Range("A2:T21").Value = $data

Write-Output "Updated $($data.GetLength(0)) rows x $($data.GetLength(1)) cols on $($ws.Name)"
